# SMLP_YR_FIN.xlsx update
# A new reporting period (column) is inserted as the new column D of the
# "SMLP" worksheet (Income Statement, Balance Sheet and Cash Flow Statement
# blocks). All existing data in columns D:K shifts right to E:L, and the new
# column D is populated with the latest period's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SMLP")

# Insert a new column before column D; this shifts the existing D:K data to
# E:L (formulas/values, shared-string references, dimension, etc. all move
# automatically).
$ws.Columns("D").Insert()

# Copy the formatting (number format / style) from the column that used to be
# D (now E) into the newly inserted column D, so the new column matches the
# look of the rest of the table (date format for the "Period Ending" rows,
# number format for the data rows). Restricted to the used data rows so we
# do not inflate the worksheet's used range down to row 1048576.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Income Statement -------------------------------------------------
$ws.Range("D7").Value  = 43465
$ws.Range("D8").Value  = 506700
$ws.Range("D9").Value  = 204500
$ws.Range("D10").Value = 302100
$ws.Range("D12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 7200
$ws.Range("D15").Value = 107100
$ws.Range("D17").Value = 371700
$ws.Range("D18").Value = 135000
$ws.Range("D20").Value = -21100
$ws.Range("D21").Value = 220600
$ws.Range("D22").Value = 60500
$ws.Range("D23").Value = 53300
$ws.Range("D24").Value = 0
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 53200
$ws.Range("D27").Value = 4300
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 0
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = 21100
$ws.Range("D33").Value = 4300
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 4300

# ---- Balance Sheet ------------------------------------------------------
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 4300
$ws.Range("D42").Value = 0
$ws.Range("D43").Value = 97900
$ws.Range("D44").Value = 0
$ws.Range("D45").Value = 4000
$ws.Range("D46").Value = 106300
$ws.Range("D47").Value = 649300
$ws.Range("D48").Value = 1963700
$ws.Range("D49").Value = 289600
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 11700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 3020600
$ws.Range("D57").Value = 38400
$ws.Range("D58").Value = 0
$ws.Range("D59").Value = 71600
$ws.Range("D60").Value = 110100
$ws.Range("D61").Value = 1257700
$ws.Range("D62").Value = 431600
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 1799300
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 318900
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 0
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 902300
$ws.Range("D77").Value = 0

# ---- Cash Flow Statement --------------------------------------------------
$ws.Range("D80").Value  = 43465
$ws.Range("D81").Value  = 4300
$ws.Range("D83").Value  = 106800
$ws.Range("D84").Value  = 0
$ws.Range("D85").Value  = 0
$ws.Range("D86").Value  = 0
$ws.Range("D87").Value  = 0
$ws.Range("D88").Value  = 0
$ws.Range("D89").Value  = 227900
$ws.Range("D91").Value  = -200600
$ws.Range("D92").Value  = 0
$ws.Range("D93").Value  = 0
$ws.Range("D94").Value  = -216300
$ws.Range("D96").Value  = -180700
$ws.Range("D97").Value  = 0
$ws.Range("D98").Value  = 0
$ws.Range("D99").Value  = 0
$ws.Range("D100").Value = -8700
$ws.Range("D101").Value = 0
$ws.Range("D102").Value = 2900
